$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B13").Value = 45147
$ws.Range("C13").Formula = "=D13-D12"
$ws.Range("D13").Value = 2721
$ws.Range("J13").Formula = "=AVERAGE(C9:C12)"
